$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet named "Sheet1" after the existing sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# Fill content for Sheet1 (order matches shared-string insertion order)
$ws2.Range("B2").Value = "CN"
$ws2.Range("C2").Value = 0
$ws2.Range("B3").Value = "Thu 2"
$ws2.Range("C3").Value = 1
$ws2.Range("B4").Value = "Thu 3"
$ws2.Range("C4").Value = 2
$ws2.Range("B5").Value = "Thu 4"
$ws2.Range("C5").Value = 3
$ws2.Range("B6").Value = "Thu 5"
$ws2.Range("C6").Value = 4
$ws2.Range("B7").Value = "Thu 6"
$ws2.Range("C7").Value = 5
$ws2.Range("B8").Value = "Thu 7"
$ws2.Range("C8").Value = 6
$ws2.Range("D8").Value = "x"
$ws2.Range("E4").Value = "y"
$ws2.Range("B11").Value = "x > y"
$ws2.Range("C12").Value = "y + (6-x) + 1"

# Column widths for D and E (target stored widths 4.42578125 / 4.5703125;
# engine snaps ColumnWidth to 1/6 increments, so 3.7 chars is the closest match)
$ws2.Columns.Item(4).ColumnWidth = 3.7
$ws2.Columns.Item(5).ColumnWidth = 3.7

# Selection on sheet1 changes from F11 to B2:B9
$ws1.Range("B2:B9").Select()

# Make Sheet1 (ws2) the active sheet, with selection E4
$ws2.Activate()
$ws2.Range("E4").Select()
